$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 41 (pushes existing rows 41-62 down to 42-63)
$ws.Rows(41).Insert()

# Populate the newly inserted row 41 with the new record's data
$ws.Range("A41").Value = 10
$ws.Range("B41").Value = "Vega Modelo de Temuco"
$ws.Range("C41").Value = "La Araucanía"
$ws.Range("D41").Value = 44795
$ws.Range("E41").Value = 9
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100108
$ws.Range("H41").Value = "Tropicales y subtropicales"
$ws.Range("I41").Value = 100108003
$ws.Range("J41").Value = "Maracuyá"
$ws.Range("K41").Value = "Sin especificar"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 25
$ws.Range("N41").Value = 36000
$ws.Range("O41").Value = 36000
$ws.Range("P41").Value = 36000
$ws.Range("Q41").Value = '$/caja 18 kilos'
$ws.Range("R41").Value = "Región de Arica y Parinacota"
$ws.Range("S41").Value = 2000
$ws.Range("T41").Value = 18
